$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. date style on col A) from the last existing row (385) down to the new rows (386-464)
$ws.Range("A385:D385").Copy()
$ws.Range("A386:D464").PasteSpecial(-4122)

$ws.Cells.Item(386,1).Value = 44460
$ws.Cells.Item(386,2).Value = 0
$ws.Cells.Item(386,3).Value = 2
$ws.Cells.Item(386,4).Value = 32.34675723758694
$ws.Cells.Item(387,1).Value = 44461
$ws.Cells.Item(387,2).Value = 0
$ws.Cells.Item(387,3).Value = 2
$ws.Cells.Item(387,4).Value = 32.34675723758694
$ws.Cells.Item(388,1).Value = 44462
$ws.Cells.Item(388,2).Value = 0
$ws.Cells.Item(388,3).Value = 1
$ws.Cells.Item(388,4).Value = 16.17337861879347
$ws.Cells.Item(389,1).Value = 44463
$ws.Cells.Item(389,2).Value = 0
$ws.Cells.Item(389,3).Value = 0
$ws.Cells.Item(389,4).Value = 0
$ws.Cells.Item(390,1).Value = 44464
$ws.Cells.Item(390,2).Value = 0
$ws.Cells.Item(390,3).Value = 0
$ws.Cells.Item(390,4).Value = 0
$ws.Cells.Item(391,1).Value = 44465
$ws.Cells.Item(391,2).Value = 0
$ws.Cells.Item(391,3).Value = 0
$ws.Cells.Item(391,4).Value = 0
$ws.Cells.Item(392,1).Value = 44466
$ws.Cells.Item(392,2).Value = 0
$ws.Cells.Item(392,3).Value = 0
$ws.Cells.Item(392,4).Value = 0
$ws.Cells.Item(393,1).Value = 44467
$ws.Cells.Item(393,2).Value = 0
$ws.Cells.Item(393,3).Value = 0
$ws.Cells.Item(393,4).Value = 0
$ws.Cells.Item(394,1).Value = 44468
$ws.Cells.Item(394,2).Value = 0
$ws.Cells.Item(394,3).Value = 0
$ws.Cells.Item(394,4).Value = 0
$ws.Cells.Item(395,1).Value = 44469
$ws.Cells.Item(395,2).Value = 0
$ws.Cells.Item(395,3).Value = 0
$ws.Cells.Item(395,4).Value = 0
$ws.Cells.Item(396,1).Value = 44470
$ws.Cells.Item(396,2).Value = 1
$ws.Cells.Item(396,3).Value = 1
$ws.Cells.Item(396,4).Value = 16.17337861879347
$ws.Cells.Item(397,1).Value = 44471
$ws.Cells.Item(397,2).Value = 0
$ws.Cells.Item(397,3).Value = 1
$ws.Cells.Item(397,4).Value = 16.17337861879347
$ws.Cells.Item(398,1).Value = 44472
$ws.Cells.Item(398,2).Value = 0
$ws.Cells.Item(398,3).Value = 1
$ws.Cells.Item(398,4).Value = 16.17337861879347
$ws.Cells.Item(399,1).Value = 44473
$ws.Cells.Item(399,2).Value = 0
$ws.Cells.Item(399,3).Value = 1
$ws.Cells.Item(399,4).Value = 16.17337861879347
$ws.Cells.Item(400,1).Value = 44474
$ws.Cells.Item(400,2).Value = 0
$ws.Cells.Item(400,3).Value = 1
$ws.Cells.Item(400,4).Value = 16.17337861879347
$ws.Cells.Item(401,1).Value = 44475
$ws.Cells.Item(401,2).Value = 0
$ws.Cells.Item(401,3).Value = 1
$ws.Cells.Item(401,4).Value = 16.17337861879347
$ws.Cells.Item(402,1).Value = 44476
$ws.Cells.Item(402,2).Value = 0
$ws.Cells.Item(402,3).Value = 1
$ws.Cells.Item(402,4).Value = 16.17337861879347
$ws.Cells.Item(403,1).Value = 44477
$ws.Cells.Item(403,2).Value = 0
$ws.Cells.Item(403,3).Value = 0
$ws.Cells.Item(403,4).Value = 0
$ws.Cells.Item(404,1).Value = 44478
$ws.Cells.Item(404,2).Value = 1
$ws.Cells.Item(404,3).Value = 1
$ws.Cells.Item(404,4).Value = 16.17337861879347
$ws.Cells.Item(405,1).Value = 44479
$ws.Cells.Item(405,2).Value = 0
$ws.Cells.Item(405,3).Value = 1
$ws.Cells.Item(405,4).Value = 16.17337861879347
$ws.Cells.Item(406,1).Value = 44480
$ws.Cells.Item(406,2).Value = 0
$ws.Cells.Item(406,3).Value = 1
$ws.Cells.Item(406,4).Value = 16.17337861879347
$ws.Cells.Item(407,1).Value = 44481
$ws.Cells.Item(407,2).Value = 0
$ws.Cells.Item(407,3).Value = 1
$ws.Cells.Item(407,4).Value = 16.17337861879347
$ws.Cells.Item(408,1).Value = 44482
$ws.Cells.Item(408,2).Value = 0
$ws.Cells.Item(408,3).Value = 1
$ws.Cells.Item(408,4).Value = 16.17337861879347
$ws.Cells.Item(409,1).Value = 44483
$ws.Cells.Item(409,2).Value = 0
$ws.Cells.Item(409,3).Value = 1
$ws.Cells.Item(409,4).Value = 16.17337861879347
$ws.Cells.Item(410,1).Value = 44484
$ws.Cells.Item(410,2).Value = 0
$ws.Cells.Item(410,3).Value = 1
$ws.Cells.Item(410,4).Value = 16.17337861879347
$ws.Cells.Item(411,1).Value = 44485
$ws.Cells.Item(411,2).Value = 0
$ws.Cells.Item(411,3).Value = 0
$ws.Cells.Item(411,4).Value = 0
$ws.Cells.Item(412,1).Value = 44486
$ws.Cells.Item(412,2).Value = 0
$ws.Cells.Item(412,3).Value = 0
$ws.Cells.Item(412,4).Value = 0
$ws.Cells.Item(413,1).Value = 44487
$ws.Cells.Item(413,2).Value = 0
$ws.Cells.Item(413,3).Value = 0
$ws.Cells.Item(413,4).Value = 0
$ws.Cells.Item(414,1).Value = 44488
$ws.Cells.Item(414,2).Value = 0
$ws.Cells.Item(414,3).Value = 0
$ws.Cells.Item(414,4).Value = 0
$ws.Cells.Item(415,1).Value = 44489
$ws.Cells.Item(415,2).Value = 0
$ws.Cells.Item(415,3).Value = 0
$ws.Cells.Item(415,4).Value = 0
$ws.Cells.Item(416,1).Value = 44490
$ws.Cells.Item(416,2).Value = 0
$ws.Cells.Item(416,3).Value = 0
$ws.Cells.Item(416,4).Value = 0
$ws.Cells.Item(417,1).Value = 44491
$ws.Cells.Item(417,2).Value = 0
$ws.Cells.Item(417,3).Value = 0
$ws.Cells.Item(417,4).Value = 0
$ws.Cells.Item(418,1).Value = 44492
$ws.Cells.Item(418,2).Value = 1
$ws.Cells.Item(418,3).Value = 1
$ws.Cells.Item(418,4).Value = 16.17337861879347
$ws.Cells.Item(419,1).Value = 44493
$ws.Cells.Item(419,2).Value = 0
$ws.Cells.Item(419,3).Value = 1
$ws.Cells.Item(419,4).Value = 16.17337861879347
$ws.Cells.Item(420,1).Value = 44494
$ws.Cells.Item(420,2).Value = 0
$ws.Cells.Item(420,3).Value = 1
$ws.Cells.Item(420,4).Value = 16.17337861879347
$ws.Cells.Item(421,1).Value = 44495
$ws.Cells.Item(421,2).Value = 0
$ws.Cells.Item(421,3).Value = 1
$ws.Cells.Item(421,4).Value = 16.17337861879347
$ws.Cells.Item(422,1).Value = 44496
$ws.Cells.Item(422,2).Value = 0
$ws.Cells.Item(422,3).Value = 1
$ws.Cells.Item(422,4).Value = 16.17337861879347
$ws.Cells.Item(423,1).Value = 44497
$ws.Cells.Item(423,2).Value = 0
$ws.Cells.Item(423,3).Value = 1
$ws.Cells.Item(423,4).Value = 16.17337861879347
$ws.Cells.Item(424,1).Value = 44498
$ws.Cells.Item(424,2).Value = 2
$ws.Cells.Item(424,3).Value = 3
$ws.Cells.Item(424,4).Value = 48.5201358563804
$ws.Cells.Item(425,1).Value = 44499
$ws.Cells.Item(425,2).Value = 1
$ws.Cells.Item(425,3).Value = 3
$ws.Cells.Item(425,4).Value = 48.5201358563804
$ws.Cells.Item(426,1).Value = 44500
$ws.Cells.Item(426,2).Value = 0
$ws.Cells.Item(426,3).Value = 3
$ws.Cells.Item(426,4).Value = 48.5201358563804
$ws.Cells.Item(427,1).Value = 44501
$ws.Cells.Item(427,2).Value = 0
$ws.Cells.Item(427,3).Value = 3
$ws.Cells.Item(427,4).Value = 48.5201358563804
$ws.Cells.Item(428,1).Value = 44502
$ws.Cells.Item(428,2).Value = 0
$ws.Cells.Item(428,3).Value = 3
$ws.Cells.Item(428,4).Value = 48.5201358563804
$ws.Cells.Item(429,1).Value = 44503
$ws.Cells.Item(429,2).Value = 0
$ws.Cells.Item(429,3).Value = 3
$ws.Cells.Item(429,4).Value = 48.5201358563804
$ws.Cells.Item(430,1).Value = 44504
$ws.Cells.Item(430,2).Value = 0
$ws.Cells.Item(430,3).Value = 3
$ws.Cells.Item(430,4).Value = 48.5201358563804
$ws.Cells.Item(431,1).Value = 44505
$ws.Cells.Item(431,2).Value = 0
$ws.Cells.Item(431,3).Value = 1
$ws.Cells.Item(431,4).Value = 16.17337861879347
$ws.Cells.Item(432,1).Value = 44506
$ws.Cells.Item(432,2).Value = 0
$ws.Cells.Item(432,3).Value = 0
$ws.Cells.Item(432,4).Value = 0
$ws.Cells.Item(433,1).Value = 44507
$ws.Cells.Item(433,2).Value = 0
$ws.Cells.Item(433,3).Value = 0
$ws.Cells.Item(433,4).Value = 0
$ws.Cells.Item(434,1).Value = 44508
$ws.Cells.Item(434,2).Value = 0
$ws.Cells.Item(434,3).Value = 0
$ws.Cells.Item(434,4).Value = 0
$ws.Cells.Item(435,1).Value = 44509
$ws.Cells.Item(435,2).Value = 0
$ws.Cells.Item(435,3).Value = 0
$ws.Cells.Item(435,4).Value = 0
$ws.Cells.Item(436,1).Value = 44510
$ws.Cells.Item(436,2).Value = 0
$ws.Cells.Item(436,3).Value = 0
$ws.Cells.Item(436,4).Value = 0
$ws.Cells.Item(437,1).Value = 44511
$ws.Cells.Item(437,2).Value = 0
$ws.Cells.Item(437,3).Value = 0
$ws.Cells.Item(437,4).Value = 0
$ws.Cells.Item(438,1).Value = 44512
$ws.Cells.Item(438,2).Value = 0
$ws.Cells.Item(438,3).Value = 0
$ws.Cells.Item(438,4).Value = 0
$ws.Cells.Item(439,1).Value = 44513
$ws.Cells.Item(439,2).Value = 1
$ws.Cells.Item(439,3).Value = 1
$ws.Cells.Item(439,4).Value = 16.17337861879347
$ws.Cells.Item(440,1).Value = 44514
$ws.Cells.Item(440,2).Value = 2
$ws.Cells.Item(440,3).Value = 3
$ws.Cells.Item(440,4).Value = 48.5201358563804
$ws.Cells.Item(441,1).Value = 44515
$ws.Cells.Item(441,2).Value = 1
$ws.Cells.Item(441,3).Value = 4
$ws.Cells.Item(441,4).Value = 64.69351447517387
$ws.Cells.Item(442,1).Value = 44516
$ws.Cells.Item(442,2).Value = 8
$ws.Cells.Item(442,3).Value = 12
$ws.Cells.Item(442,4).Value = 194.0805434255216
$ws.Cells.Item(443,1).Value = 44517
$ws.Cells.Item(443,2).Value = 0
$ws.Cells.Item(443,3).Value = 12
$ws.Cells.Item(443,4).Value = 194.0805434255216
$ws.Cells.Item(444,1).Value = 44518
$ws.Cells.Item(444,2).Value = 0
$ws.Cells.Item(444,3).Value = 12
$ws.Cells.Item(444,4).Value = 194.0805434255216
$ws.Cells.Item(445,1).Value = 44519
$ws.Cells.Item(445,2).Value = 1
$ws.Cells.Item(445,3).Value = 13
$ws.Cells.Item(445,4).Value = 210.2539220443151
$ws.Cells.Item(446,1).Value = 44520
$ws.Cells.Item(446,2).Value = 0
$ws.Cells.Item(446,3).Value = 12
$ws.Cells.Item(446,4).Value = 194.0805434255216
$ws.Cells.Item(447,1).Value = 44521
$ws.Cells.Item(447,2).Value = 0
$ws.Cells.Item(447,3).Value = 10
$ws.Cells.Item(447,4).Value = 161.7337861879347
$ws.Cells.Item(448,1).Value = 44522
$ws.Cells.Item(448,2).Value = 0
$ws.Cells.Item(448,3).Value = 9
$ws.Cells.Item(448,4).Value = 145.5604075691412
$ws.Cells.Item(449,1).Value = 44523
$ws.Cells.Item(449,2).Value = 0
$ws.Cells.Item(449,3).Value = 1
$ws.Cells.Item(449,4).Value = 16.17337861879347
$ws.Cells.Item(450,1).Value = 44524
$ws.Cells.Item(450,2).Value = 6
$ws.Cells.Item(450,3).Value = 7
$ws.Cells.Item(450,4).Value = 113.2136503315543
$ws.Cells.Item(451,1).Value = 44525
$ws.Cells.Item(451,2).Value = 0
$ws.Cells.Item(451,3).Value = 7
$ws.Cells.Item(451,4).Value = 113.2136503315543
$ws.Cells.Item(452,1).Value = 44526
$ws.Cells.Item(452,2).Value = 3
$ws.Cells.Item(452,3).Value = 9
$ws.Cells.Item(452,4).Value = 145.5604075691412
$ws.Cells.Item(453,1).Value = 44527
$ws.Cells.Item(453,2).Value = 2
$ws.Cells.Item(453,3).Value = 11
$ws.Cells.Item(453,4).Value = 177.9071648067281
$ws.Cells.Item(454,1).Value = 44528
$ws.Cells.Item(454,2).Value = 0
$ws.Cells.Item(454,3).Value = 11
$ws.Cells.Item(454,4).Value = 177.9071648067281
$ws.Cells.Item(455,1).Value = 44529
$ws.Cells.Item(455,2).Value = 1
$ws.Cells.Item(455,3).Value = 12
$ws.Cells.Item(455,4).Value = 194.0805434255216
$ws.Cells.Item(456,1).Value = 44530
$ws.Cells.Item(456,2).Value = 1
$ws.Cells.Item(456,3).Value = 13
$ws.Cells.Item(456,4).Value = 210.2539220443151
$ws.Cells.Item(457,1).Value = 44531
$ws.Cells.Item(457,2).Value = 0
$ws.Cells.Item(457,3).Value = 7
$ws.Cells.Item(457,4).Value = 113.2136503315543
$ws.Cells.Item(458,1).Value = 44532
$ws.Cells.Item(458,2).Value = 3
$ws.Cells.Item(458,3).Value = 10
$ws.Cells.Item(458,4).Value = 161.7337861879347
$ws.Cells.Item(459,1).Value = 44533
$ws.Cells.Item(459,2).Value = 9
$ws.Cells.Item(459,3).Value = 16
$ws.Cells.Item(459,4).Value = 258.7740579006955
$ws.Cells.Item(460,1).Value = 44534
$ws.Cells.Item(460,2).Value = 1
$ws.Cells.Item(460,3).Value = 15
$ws.Cells.Item(460,4).Value = 242.600679281902
$ws.Cells.Item(461,1).Value = 44535
$ws.Cells.Item(461,2).Value = 1
$ws.Cells.Item(461,3).Value = 16
$ws.Cells.Item(461,4).Value = 258.7740579006955
$ws.Cells.Item(462,1).Value = 44536
$ws.Cells.Item(462,2).Value = 8
$ws.Cells.Item(462,3).Value = 23
$ws.Cells.Item(462,4).Value = 371.9877082322497
$ws.Cells.Item(463,1).Value = 44537
$ws.Cells.Item(463,2).Value = 0
$ws.Cells.Item(463,3).Value = 22
$ws.Cells.Item(463,4).Value = 355.8143296134562
$ws.Cells.Item(464,1).Value = 44538
$ws.Cells.Item(464,2).Value = 0
$ws.Cells.Item(464,3).Value = 22
$ws.Cells.Item(464,4).Value = 355.8143296134562
